# Auto-generated edit script: updates Leve profit-calculation values
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match
# refreshed market-board prices from the scheduled Chocobo_Profits runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 3150  # H76
$ws.Cells.Item(76, 9).Value = 3150  # I76
$ws.Cells.Item(76, 11).Value = 3150  # K76
$ws.Cells.Item(76, 13).Value = -2835  # M76
$ws.Cells.Item(79, 8).Value = 3150  # H79
$ws.Cells.Item(79, 9).Value = 3150  # I79
$ws.Cells.Item(79, 11).Value = 3150  # K79
$ws.Cells.Item(79, 13).Value = -2058  # M79
$ws.Cells.Item(129, 8).Value = 1053.7037  # H129
$ws.Cells.Item(129, 9).Value = 333.33334  # I129
$ws.Cells.Item(129, 10).Value = 1096.0785  # J129
$ws.Cells.Item(129, 11).Value = 1000.00002  # K129
$ws.Cells.Item(129, 12).Value = 3288.2355  # L129
$ws.Cells.Item(129, 13).Value = 3999.99998  # M129
$ws.Cells.Item(129, 14).Value = -13288.2355  # N129
$ws.Cells.Item(132, 8).Value = 97742.69  # H132
$ws.Cells.Item(132, 9).Value = 107005.08  # I132
$ws.Cells.Item(132, 11).Value = 321015.24  # K132
$ws.Cells.Item(132, 13).Value = -318485.24  # M132
$ws.Cells.Item(139, 8).Value = 41375.715  # H139
$ws.Cells.Item(139, 10).Value = 41375.715  # J139
$ws.Cells.Item(139, 12).Value = 41375.715  # L139
$ws.Cells.Item(139, 14).Value = -51655.715  # N139
$ws.Cells.Item(141, 8).Value = 1907.875  # H141
$ws.Cells.Item(141, 9).Value = 1725.7587  # I141
$ws.Cells.Item(141, 10).Value = 3668.3333  # J141
$ws.Cells.Item(141, 11).Value = 5177.2761  # K141
$ws.Cells.Item(141, 12).Value = 11004.9999  # L141
$ws.Cells.Item(141, 13).Value = 2.723899999999958  # M141
$ws.Cells.Item(141, 14).Value = -21364.9999  # N141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6588.5063  # H32
$ws.Cells.Item(32, 9).Value = 3788.647  # I32
$ws.Cells.Item(32, 10).Value = 12080.538  # J32
$ws.Cells.Item(32, 11).Value = 3788.647  # K32
$ws.Cells.Item(32, 12).Value = 12080.538  # L32
$ws.Cells.Item(32, 13).Value = -3501.647  # M32
$ws.Cells.Item(32, 14).Value = -12654.538  # N32
$ws.Cells.Item(74, 8).Value = 8217  # H74
$ws.Cells.Item(74, 9).Value = 8945  # I74
$ws.Cells.Item(74, 11).Value = 8945  # K74
$ws.Cells.Item(74, 13).Value = -8071  # M74
$ws.Cells.Item(77, 8).Value = 8217  # H77
$ws.Cells.Item(77, 9).Value = 8945  # I77
$ws.Cells.Item(77, 11).Value = 44725  # K77
$ws.Cells.Item(77, 13).Value = -40357  # M77
$ws.Cells.Item(88, 8).Value = 33334832  # H88
$ws.Cells.Item(88, 9).Value = 33334832  # I88
$ws.Cells.Item(88, 11).Value = 33334832  # K88
$ws.Cells.Item(88, 13).Value = -33334426  # M88
$ws.Cells.Item(91, 8).Value = 33334832  # H91
$ws.Cells.Item(91, 9).Value = 33334832  # I91
$ws.Cells.Item(91, 11).Value = 33334832  # K91
$ws.Cells.Item(91, 13).Value = -33333428  # M91
$ws.Cells.Item(132, 8).Value = 3071.4707  # H132
$ws.Cells.Item(132, 9).Value = 2059.6667  # I132
$ws.Cells.Item(132, 10).Value = 5499.8  # J132
$ws.Cells.Item(132, 11).Value = 6179.000100000001  # K132
$ws.Cells.Item(132, 12).Value = 16499.4  # L132
$ws.Cells.Item(132, 13).Value = -3649.000100000001  # M132
$ws.Cells.Item(132, 14).Value = -21559.4  # N132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(63, 8).Value = 55000  # H63
$ws.Cells.Item(63, 10).Value = 55000  # J63
$ws.Cells.Item(63, 12).Value = 55000  # L63
$ws.Cells.Item(63, 14).Value = -56372  # N63
$ws.Cells.Item(66, 8).Value = 55000  # H66
$ws.Cells.Item(66, 10).Value = 55000  # J66
$ws.Cells.Item(66, 12).Value = 165000  # L66
$ws.Cells.Item(66, 14).Value = -171864  # N66
$ws.Cells.Item(134, 8).Value = 2843.3962  # H134
$ws.Cells.Item(134, 9).Value = 1629.2069  # I134
$ws.Cells.Item(134, 10).Value = 4310.5415  # J134
$ws.Cells.Item(134, 11).Value = 4887.620699999999  # K134
$ws.Cells.Item(134, 12).Value = 12931.6245  # L134
$ws.Cells.Item(134, 13).Value = -2352.620699999999  # M134
$ws.Cells.Item(134, 14).Value = -18001.6245  # N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(9, 8).Value = 20950  # H9
$ws.Cells.Item(9, 10).Value = 20950  # J9
$ws.Cells.Item(9, 12).Value = 20950  # L9
$ws.Cells.Item(9, 14).Value = -21286  # N9
$ws.Cells.Item(62, 8).Value = 125004750  # H62
$ws.Cells.Item(62, 9).Value = 125004750  # I62
$ws.Cells.Item(62, 11).Value = 125004750  # K62
$ws.Cells.Item(62, 13).Value = -125004126  # M62
$ws.Cells.Item(65, 8).Value = 125004750  # H65
$ws.Cells.Item(65, 9).Value = 125004750  # I65
$ws.Cells.Item(65, 11).Value = 625023750  # K65
$ws.Cells.Item(65, 13).Value = -625020630  # M65
$ws.Cells.Item(87, 8).Value = 23714.285  # H87
$ws.Cells.Item(87, 10).Value = 23714.285  # J87
$ws.Cells.Item(87, 12).Value = 23714.285  # L87
$ws.Cells.Item(87, 14).Value = -26086.285  # N87
$ws.Cells.Item(90, 8).Value = 23714.285  # H90
$ws.Cells.Item(90, 10).Value = 23714.285  # J90
$ws.Cells.Item(90, 12).Value = 71142.855  # L90
$ws.Cells.Item(90, 14).Value = -82998.855  # N90
$ws.Cells.Item(132, 8).Value = 3073.2856  # H132
$ws.Cells.Item(132, 9).Value = 1398.4117  # I132
$ws.Cells.Item(132, 10).Value = 4655.1113  # J132
$ws.Cells.Item(132, 11).Value = 4195.2351  # K132
$ws.Cells.Item(132, 12).Value = 13965.3339  # L132
$ws.Cells.Item(132, 13).Value = -1665.2351  # M132
$ws.Cells.Item(132, 14).Value = -19025.3339  # N132
$ws.Cells.Item(134, 8).Value = 7821.1665  # H134
$ws.Cells.Item(134, 9).Value = 7718.7334  # I134
$ws.Cells.Item(134, 11).Value = 23156.2002  # K134
$ws.Cells.Item(134, 13).Value = -20621.2002  # M134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1025.0725  # H131
$ws.Cells.Item(131, 10).Value = 915.0769  # J131
$ws.Cells.Item(131, 12).Value = 2745.2307  # L131
$ws.Cells.Item(131, 14).Value = -12825.2307  # N131
$ws.Cells.Item(137, 8).Value = 3552.0715  # H137
$ws.Cells.Item(137, 9).Value = 1505  # I137
$ws.Cells.Item(137, 10).Value = 5087.375  # J137
$ws.Cells.Item(137, 11).Value = 4515  # K137
$ws.Cells.Item(137, 12).Value = 15262.125  # L137
$ws.Cells.Item(137, 13).Value = 585  # M137
$ws.Cells.Item(137, 14).Value = -25462.125  # N137
$ws.Cells.Item(140, 8).Value = 20526.172  # H140
$ws.Cells.Item(140, 9).Value = 45094.082  # I140
$ws.Cells.Item(140, 10).Value = 3184.1177  # J140
$ws.Cells.Item(140, 11).Value = 135282.246  # K140
$ws.Cells.Item(140, 12).Value = 9552.3531  # L140
$ws.Cells.Item(140, 13).Value = -130102.246  # M140
$ws.Cells.Item(140, 14).Value = -19912.3531  # N140

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6783.1377  # H70
$ws.Cells.Item(70, 9).Value = 5860.048  # I70
$ws.Cells.Item(70, 10).Value = 9206.25  # J70
$ws.Cells.Item(70, 11).Value = 5860.048  # K70
$ws.Cells.Item(70, 12).Value = 9206.25  # L70
$ws.Cells.Item(70, 13).Value = -5590.048  # M70
$ws.Cells.Item(70, 14).Value = -9746.25  # N70
$ws.Cells.Item(73, 8).Value = 6783.1377  # H73
$ws.Cells.Item(73, 9).Value = 5860.048  # I73
$ws.Cells.Item(73, 10).Value = 9206.25  # J73
$ws.Cells.Item(73, 11).Value = 5860.048  # K73
$ws.Cells.Item(73, 12).Value = 9206.25  # L73
$ws.Cells.Item(73, 13).Value = -4924.048  # M73
$ws.Cells.Item(73, 14).Value = -11078.25  # N73

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 1927.9546  # H93
$ws.Cells.Item(93, 9).Value = 1131.6154  # I93
$ws.Cells.Item(93, 11).Value = 1131.6154  # K93
$ws.Cells.Item(93, 13).Value = 116.3846000000001  # M93

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 5495375  # H81
$ws.Cells.Item(81, 9).Value = 7937457  # I81
$ws.Cells.Item(81, 10).Value = 691.25  # J81
$ws.Cells.Item(81, 11).Value = 15874914  # K81
$ws.Cells.Item(81, 12).Value = 1382.5  # L81
$ws.Cells.Item(81, 13).Value = -15873853  # M81
$ws.Cells.Item(81, 14).Value = -3504.5  # N81
$ws.Cells.Item(84, 8).Value = 5495375  # H84
$ws.Cells.Item(84, 9).Value = 7937457  # I84
$ws.Cells.Item(84, 10).Value = 691.25  # J84
$ws.Cells.Item(84, 11).Value = 79374570  # K84
$ws.Cells.Item(84, 12).Value = 6912.5  # L84
$ws.Cells.Item(84, 13).Value = -79369266  # M84
$ws.Cells.Item(84, 14).Value = -17520.5  # N84
$ws.Cells.Item(132, 8).Value = 12822870  # H132
$ws.Cells.Item(132, 9).Value = 862.36365  # I132
$ws.Cells.Item(132, 10).Value = 22225676  # J132
$ws.Cells.Item(132, 11).Value = 2587.09095  # K132
$ws.Cells.Item(132, 12).Value = 66677028  # L132
$ws.Cells.Item(132, 13).Value = -57.09094999999979  # M132
$ws.Cells.Item(132, 14).Value = -66682088  # N132
